# Purchase entry UI changes
$wb = $excel.ActiveWorkbook

$gatepass = $wb.Worksheets.Item("Gatepass")
$purchaseEntry = $wb.Worksheets.Item("PurchaseEntry")

# Update the PurchaseEntry invoice/medicine values.
# Set B4 first, then A4, A2, A3 so the shared-string table is rebuilt
# in the same order the source workbook expects.
$purchaseEntry.Range("B4").Value = "GOLD COAT LOTION 120ML"
$purchaseEntry.Range("A4").Value = "ne31"
$purchaseEntry.Range("A2").Value = "ne12"
$purchaseEntry.Range("A3").Value = "ne41"

# Update the selection on the (now inactive) Gatepass sheet.
$gatepass.Activate()
[void]$gatepass.Range("B4").Select()

# Make PurchaseEntry the active tab/sheet with the new selection.
$purchaseEntry.Activate()
[void]$purchaseEntry.Range("B6").Select()
